{"js": "// Update the date line and the 5x20 grid of addition/subtraction problems.\n// The mapping below reproduces the unified diff exactly: every <w:t> run\n// keeps its position in the document, only its text changes.\n\nconst newDate = \"2025-10-30 Thursday\";\n\nconst newValues = [\n  [\"99-89=\", \"63+0=\", \"51-50=\", \"10+17=\", \"86-5=\"],\n  [\"8+5=\", \"43-33=\", \"72-8=\", \"42+33=\", \"76-14=\"],\n  [\"62-62=\", \"63-8=\", \"85-1=\", \"75-33=\", \"59+8=\"],\n  [\"89-34=\", \"2+14=\", \"28+1=\", \"21-17=\", \"51-4=\"],\n  [\"75-44=\", \"89-30=\", \"63+23=\", \"50-23=\", \"39+43=\"],\n  [\"70+15=\", \"56-43=\", \"64-39=\", \"56-12=\", \"65-63=\"],\n  [\"79-6=\", \"33+20=\", \"5+38=\", \"97-93=\", \"76-37=\"],\n  [\"90-23=\", \"71-26=\", \"35-10=\", \"0+2=\", \"96-26=\"],\n  [\"94-62=\", \"83+7=\", \"52-1=\", \"20+8=\", \"26+53=\"],\n  [\"76-23=\", \"15+67=\", \"69+16=\", \"26+28=\", \"62+20=\"],\n  [\"18+58=\", \"56-14=\", \"47-41=\", \"97-67=\", \"64-3=\"],\n  [\"54+36=\", \"63+6=\", \"76-17=\", \"85-64=\", \"43+8=\"],\n  [\"36-1=\", \"66-49=\", \"18+66=\", \"35+53=\", \"24+59=\"],\n  [\"28+61=\", \"10+4=\", \"71-45=\", \"31-17=\", \"43-6=\"],\n  [\"2+8=\", \"74-10=\", \"8+21=\", \"72-4=\", \"99-71=\"],\n  [\"67-20=\", \"93-72=\", \"99-81=\", \"97-2=\", \"25-19=\"],\n  [\"9+47=\", \"54+17=\", \"77+17=\", \"42-35=\", \"70+13=\"],\n  [\"27-15=\", \"79+6=\", \"39+23=\", \"42+49=\", \"33+34=\"],\n  [\"20+34=\", \"84-36=\", \"5+50=\", \"63-3=\", \"80-78=\"],\n  [\"45-34=\", \"84-1=\", \"72+4=\", \"97-8=\", \"93-87=\"]\n];\n\n// 1) Update the heading paragraph with the date (first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nconst firstParagraph = paragraphs.getFirst();\nfirstParagraph.load(\"text\");\nawait context.sync();\nif (firstParagraph.text !== newDate) {\n  firstParagraph.insertText(newDate, \"Replace\");\n}\n\n// 2) Update every cell in the table of math problems, preserving formatting.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newValues;\n\nawait context.sync();\n", "ps1": "# Update the date line and the 5x20 grid of addition/subtraction problems.\n# The mapping below reproduces the unified diff exactly: every text run keeps\n# its position in the document/table, only its text changes.\n\n$d = $word.ActiveDocument\n\n# 1) Update the heading paragraph with the date (first paragraph in the body).\n$d.Paragraphs.First.Range.Text = \"2025-10-30 Thursday\"\n\n# 2) Update every cell in the table of math problems, preserving formatting.\n$newValues = @(\n    @(\"99-89=\", \"63+0=\", \"51-50=\", \"10+17=\", \"86-5=\"),\n    @(\"8+5=\", \"43-33=\", \"72-8=\", \"42+33=\", \"76-14=\"),\n    @(\"62-62=\", \"63-8=\", \"85-1=\", \"75-33=\", \"59+8=\"),\n    @(\"89-34=\", \"2+14=\", \"28+1=\", \"21-17=\", \"51-4=\"),\n    @(\"75-44=\", \"89-30=\", \"63+23=\", \"50-23=\", \"39+43=\"),\n    @(\"70+15=\", \"56-43=\", \"64-39=\", \"56-12=\", \"65-63=\"),\n    @(\"79-6=\", \"33+20=\", \"5+38=\", \"97-93=\", \"76-37=\"),\n    @(\"90-23=\", \"71-26=\", \"35-10=\", \"0+2=\", \"96-26=\"),\n    @(\"94-62=\", \"83+7=\", \"52-1=\", \"20+8=\", \"26+53=\"),\n    @(\"76-23=\", \"15+67=\", \"69+16=\", \"26+28=\", \"62+20=\"),\n    @(\"18+58=\", \"56-14=\", \"47-41=\", \"97-67=\", \"64-3=\"),\n    @(\"54+36=\", \"63+6=\", \"76-17=\", \"85-64=\", \"43+8=\"),\n    @(\"36-1=\", \"66-49=\", \"18+66=\", \"35+53=\", \"24+59=\"),\n    @(\"28+61=\", \"10+4=\", \"71-45=\", \"31-17=\", \"43-6=\"),\n    @(\"2+8=\", \"74-10=\", \"8+21=\", \"72-4=\", \"99-71=\"),\n    @(\"67-20=\", \"93-72=\", \"99-81=\", \"97-2=\", \"25-19=\"),\n    @(\"9+47=\", \"54+17=\", \"77+17=\", \"42-35=\", \"70+13=\"),\n    @(\"27-15=\", \"79+6=\", \"39+23=\", \"42+49=\", \"33+34=\"),\n    @(\"20+34=\", \"84-36=\", \"5+50=\", \"63-3=\", \"80-78=\"),\n    @(\"45-34=\", \"84-1=\", \"72+4=\", \"97-8=\", \"93-87=\"),\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
